# Rename worksheets
$wb = $excel.ActiveWorkbook
$wb.Worksheets.Item(1).Name = "N"
$wb.Worksheets.Item(2).Name = "2N"
$wb.Worksheets.Item(3).Name = "3N"

# Sheet "N" (was "Simulation 1"): reset the scrolled view back to the top-left
# cell (removes topLeftCell="A34"), keep the existing selection.
$wsN = $wb.Worksheets.Item(1)
$wsN.Activate()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1

# Sheet "3N" (was "Simulation 3"): change the zoom level from 25% to 65%,
# and keep it the active/selected tab when done (matches activeTab=2).
$ws3N = $wb.Worksheets.Item(3)
$ws3N.Activate()
$excel.ActiveWindow.Zoom = 65
